$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing column D ("batsman"), shifting
# batsman/totalRuns/totalBalls/total4s/total6s/sr from D:I to F:K, to make
# room for the new ownTeam / oppTeam columns.
$ws.Range("D1:E1").EntireColumn.Insert()

# Header row
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# Row 2 data
$ws.Range("D2").Value = "Mumbai Indians"
$ws.Range("E2").Value = "Sunrisers Hyderabad"

# Row 3 data
$ws.Range("D3").Value = "Mumbai Indians"
$ws.Range("E3").Value = "Kings XI Punjab"
